# Remove the two unused variables "GoalTime" and "GoalDistance" (rows 19
# and 20) from Sheet1, shifting everything below them up by two rows, and
# update the sheet's selection to match the post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 19 and 20 (GoalTime, GoalDistance) entirely -- remaining rows
# shift up to fill the gap.
$ws.Rows("19:20").Delete()

# Update the visible selection on the sheet to match the saved state.
$ws.Range("A19:XFD20").Select()
